# The dataset sheet had a blank first row (data started at row 2, A2:B104).
# This edit removes that leading blank row, which shifts every data row up
# by one (A2->A1 ... A104 disappears, new used range A1:B103).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(1).Delete()

# Re-apply the (pre-existing) sort so the sheet's remembered sort range
# follows the data up by one row as well.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C103"))
$ws.Sort.SetRange($ws.Range("A1:C103"))
$ws.Sort.Header = -4163
$ws.Sort.Apply()

# After deleting the row, Excel leaves the freed row selected.
$ws.Rows(1).Select()
